$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Search users. (Needs back-end)" -> "(DONE) Search users."
#    split into three runs: "(DONE) " | "Search users" | "."
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Search users. (Needs back-end)`r") {
        $rng = $p.Range
        $rng.Find.Execute("Search users. (Needs back-end)", $true, $false, $false, $false, $false, $true, 1, $false, "(DONE) Search users.", 2)

        $s = $p.Range.Start

        # split off the leading "(DONE) " (7 chars) into its own run
        $r1 = $d.Range($s, $s + 7)
        $r1.Font.Bold = 1
        $r1.Font.Bold = 0

        # split off the trailing "." (right after "Search users", 19 chars in)
        $s2 = $p.Range.Start
        $r2 = $d.Range($s2, $s2 + 19)
        $r2.Font.Bold = 1
        $r2.Font.Bold = 0

        break
    }
}

# ------------------------------------------------------------------
# 2) Move the Word-managed "_GoBack" bookmark out of the "Logout option ..."
#    paragraph and merge its two runs back into one; re-insert the
#    bookmark (collapsed) inside the "PUT users/about - " paragraph,
#    right before "update logged in user's about."
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Logout option witho" + "ut manually changing URL.`r") {
        $p.Range.Find.Execute("Logout option without manually changing URL.", $true, $false, $false, $false, $false, $true, 1, $false, "Logout option without manually changing URL.", 2)
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*PUT users/about - update logged in user's about.*") {
        $s = $p.Range.Start
        $target = $s + 25
        $bmRange = $d.Range($target, $target)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}

# ------------------------------------------------------------------
# 3) Prepend "(DONE) " as its own run in front of the
#    "GET users/:username?partial=..." paragraph.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "GET users/:username?partial*") {
        $s = $p.Range.Start
        $ins = $d.Range($s, $s)
        $ins.InsertBefore("(DONE) ")
        break
    }
}
